$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Round Q3 and R3 to nearest integer
$ws.Range("Q3").Value2 = [math]::Round($ws.Range("Q3").Value2, 0)
$ws.Range("R3").Value2 = [math]::Round($ws.Range("R3").Value2, 0)

# Clear the Starttid (Z3) and Sluttid (AB3) cells entirely
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
